$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Simulation" group (group 1) gains a 5th loop. Its two rows (Motor /
# Premotor) are inserted right before the old row 10, pushing every row
# below down by two. Insert two blank rows at 10:11 for this purpose.
$ws.Rows("10:11").Insert()

# New row 10: Simulation, loop 5, Motor
$ws.Cells.Item(10, 1).Value = 1
$ws.Cells.Item(10, 2).Value = "Simulation"
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = "Motor"
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 347
$ws.Cells.Item(10, 7).Formula = "=F10-E10"

# New row 11: Simulation, loop 5, Premotor
$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = "Simulation"
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = "Premotor"
$ws.Cells.Item(11, 5).Formula = "=F10"
$ws.Cells.Item(11, 6).Value = 410
$ws.Cells.Item(11, 7).Formula = "=F11-E11"

# Rename the two "Exoskeleton" group names to "Experiment" everywhere they
# are used (rows 12-23 after the insert above).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
for ($r = 12; $r -le $lastRow; $r++) {
  $name = $ws.Cells.Item($r, 2).Value2
  if ($name -eq "Exoskeleton (opposition sensor at actuator)") {
    $ws.Cells.Item($r, 2).Value = "Experiment (opposition sensor at actuator)"
  } elseif ($name -eq "Exoskeleton (opposition sensor between fingers)") {
    $ws.Cells.Item($r, 2).Value = "Experiment (opposition sensor between fingers)"
  }
}

# Match the saved view/selection state.
$ws.Range("F12").Select()
